$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: department (replace school name with category)
$ws.Range("C2").Value = "Air-Conditioning"
$ws.Range("C3").Value = "Air-Conditioning"
$ws.Range("C4").Value = "Air-Conditioning"
$ws.Range("C5").Value = "Automotive"
$ws.Range("C6").Value = "Automotive"
$ws.Range("C7").Value = "Automotive"
$ws.Range("C8").Value = "Automotive"
$ws.Range("C9").Value = "Packages"
$ws.Range("C10").Value = "Packages"
$ws.Range("C11").Value = "Packages"

# Column M: location - normalise "NSW/QLD (Currently not accepting enrolments)" -> "NSW/QLD"
$ws.Range("M2").Value = "NSW/QLD"
$ws.Range("M3").Value = "NSW/QLD"
$ws.Range("M4").Value = "NSW/QLD"
$ws.Range("M9").Value = "NSW/QLD"

# Column N: locationDetail - new note for rows that are no longer accepting enrolments
$ws.Range("N2").Value = "Currently not accepting enrolments"
$ws.Range("N3").Value = "Currently not accepting enrolments"
$ws.Range("N4").Value = "Currently not accepting enrolments"
$ws.Range("N9").Value = "Currently not accepting enrolments"

# Column R: promotionValidity - clear the expired promotion text (keep formatting)
$ws.Range("R2:R11").Value = ""

Write-Output "done"
